# NYPD CompStat weekly workbook refresh:
# - bump the report "Volume ... Number" counter
# - roll the "Report Covering the Week" date range forward
# - update this week's crime-complaint statistics table (rows 15-33)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
# "Volume 32   Number  46" -> "Volume 32   Number  47"
$cell = $ws.Range("A8")
$full = $cell.Characters(1, 200)
$t = $full.Text
$t = $t.Replace("Number  46", "Number  47")
$cell.Value = $t

# "Report Covering the Week  11/10/2025  Through  11/16/2025"
# -> "Report Covering the Week  11/17/2025  Through  11/23/2025"
$cell2 = $ws.Range("C9")
$full2 = $cell2.Characters(1, 200)
$t2 = $full2.Text
$t2 = $t2.Replace("11/10/2025", "11/17/2025").Replace("11/16/2025", "11/23/2025")
$cell2.Value = $t2

# --- Row 15 (Rape) ---------------------------------------------------------
$ws.Range("C15").Value = 1
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = 20
$ws.Range("I15").Value = 57
$ws.Range("K15").Value = 83.870967741935
$ws.Range("L15").Value = 46.153846153846
$ws.Range("M15").Value = 72.727272727272
$ws.Range("N15").Value = -14.925373134328

# --- Row 16 (Robbery) -------------------------------------------------------
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 25
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = -7.407407407407
$ws.Range("J16").Value = 440
$ws.Range("K16").Value = -4.545454545454
$ws.Range("L16").Value = -13.043478260869
$ws.Range("M16").Value = -10.638297872340
$ws.Range("N16").Value = -77.443609022556

# --- Row 17 (Fel. Assault) ---------------------------------------------------
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 25
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 60
$ws.Range("G17").Value = 67
$ws.Range("H17").Value = -10.447761194029
$ws.Range("I17").Value = 719
$ws.Range("J17").Value = 713
$ws.Range("K17").Value = 0.841514726507
$ws.Range("L17").Value = -5.394736842105
$ws.Range("M17").Value = 55.627705627705
$ws.Range("N17").Value = -15.906432748538

# --- Row 18 (Burglary) -------------------------------------------------------
$ws.Range("C18").Value = 6
$ws.Range("E18").Value = 20
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = 5.263157894736
$ws.Range("I18").Value = 245
$ws.Range("J18").Value = 279
$ws.Range("K18").Value = -12.186379928315
$ws.Range("L18").Value = -16.095890410958
$ws.Range("M18").Value = -29.190751445086
$ws.Range("N18").Value = -85.494375370041

# --- Row 19 (Gr. Larceny) -----------------------------------------------------
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 24
$ws.Range("E19").Value = -45.833333333333
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 90
$ws.Range("H19").Value = -34.444444444444
$ws.Range("I19").Value = 983
$ws.Range("J19").Value = 955
$ws.Range("K19").Value = 2.931937172774
$ws.Range("L19").Value = 18.576598311218
$ws.Range("M19").Value = 90.873786407767
$ws.Range("N19").Value = 38.841807909604

# --- Row 20 (G.L.A.) ----------------------------------------------------------
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 12
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 49
$ws.Range("H20").Value = -53.061224489795
$ws.Range("I20").Value = 425
$ws.Range("J20").Value = 490
$ws.Range("K20").Value = -13.265306122449
$ws.Range("L20").Value = -27.966101694915
$ws.Range("M20").Value = 93.181818181818
$ws.Range("N20").Value = -76.558190843905

# --- Row 21 (TOTAL) -----------------------------------------------------------
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 68
$ws.Range("E21").Value = -42.647058823529
$ws.Range("F21").Value = 196
$ws.Range("G21").Value = 258
$ws.Range("H21").Value = -24.031007751938
$ws.Range("I21").Value = 2860
$ws.Range("J21").Value = 2920
$ws.Range("K21").Value = -2.054794520547
$ws.Range("L21").Value = -4.793608521970
$ws.Range("M21").Value = 38.902379796017
$ws.Range("N21").Value = -59.461374911410

# --- Row 22 (Transit) ---------------------------------------------------------
$ws.Range("G22").Value = 3

# --- Row 23 (Housing) ---------------------------------------------------------
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 7
$ws.Range("E23").Value = -71.428571428571
$ws.Range("F23").Value = 18
$ws.Range("G23").Value = 31
$ws.Range("H23").Value = -41.935483870967
$ws.Range("I23").Value = 238
$ws.Range("J23").Value = 269
$ws.Range("K23").Value = -11.524163568773
$ws.Range("L23").Value = -15
$ws.Range("M23").Value = 17.821782178217

# --- Row 24 (Petit Larceny) -----------------------------------------------------
$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 54
$ws.Range("E24").Value = -27.777777777777
$ws.Range("F24").Value = 156
$ws.Range("G24").Value = 163
$ws.Range("H24").Value = -4.294478527607
$ws.Range("I24").Value = 1849
$ws.Range("J24").Value = 1707
$ws.Range("K24").Value = 8.318687756297
$ws.Range("L24").Value = 6.632064590542
$ws.Range("M24").Value = 36.457564575645

# --- Row 25 (Retail Theft) -------------------------------------------------------
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 29
$ws.Range("E25").Value = -65.517241379310
$ws.Range("F25").Value = 46
$ws.Range("G25").Value = 60
$ws.Range("H25").Value = -23.333333333333
$ws.Range("I25").Value = 572
$ws.Range("J25").Value = 673
$ws.Range("K25").Value = -15.007429420505
$ws.Range("L25").Value = -19.549929676512

# --- Row 26 (Misd. Assault) -------------------------------------------------------
$ws.Range("C26").Value = 28
$ws.Range("D26").Value = 20
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 68
$ws.Range("H26").Value = -9.333333333333
$ws.Range("I26").Value = 977
$ws.Range("J26").Value = 1059
$ws.Range("K26").Value = -7.743153918791
$ws.Range("L26").Value = -2.3
$ws.Range("M26").Value = -30.462633451957

# --- Row 27 (UCR Rape*) -------------------------------------------------------
$ws.Range("C27").Value = 1
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -14.285714285714
$ws.Range("I27").Value = 68
$ws.Range("K27").Value = 41.666666666666
$ws.Range("L27").Value = 19.298245614035

# --- Row 28 (Other Sex Crimes) -------------------------------------------------------
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -75
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = -9.090909090909
$ws.Range("I28").Value = 70
$ws.Range("J28").Value = 103
$ws.Range("K28").Value = -32.038834951456
$ws.Range("L28").Value = -28.571428571428

# --- Row 29 (Shooting Vic.) -------------------------------------------------------
# C29 goes from a numeric 1 to a blank-style text "0" (same look as D29)
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("D29").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("L29").Value = -12.5
$ws.Range("M29").Value = -27.083333333333
$ws.Range("N29").Value = -75.177304964539

# --- Row 30 (Shooting Inc.) -------------------------------------------------------
# C30 goes from a numeric 1 to a blank-style text "0" (same look as D30)
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("D30").Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("L30").Value = -11.764705882352
$ws.Range("M30").Value = -25
$ws.Range("N30").Value = -76.377952755905

# --- Row 33 (Traffic Fatalities) -------------------------------------------------------
# D33 goes from numeric 2 to text "0" (same look as C33); E33 from numeric -100 to text "***.*"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0"
$ws.Range("C33").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4122)

$ws.Range("E33").Value = "***.*"
$ws.Range("F33").Copy() | Out-Null
$ws.Range("E33").PasteSpecial(-4122)

$ws.Range("G33").Value = 2
$ws.Range("L33").Value = -57.142857142857
